$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Value updates (new "PS-1.x" sub-task rows + DRV-2/3/4 tasks shifted down as
# the two new PS-1.9 / PS-1.10 user-story rows were inserted; gesture-related
# tasks DRV-2.x .. DRV-4.2 marked complete / in progress; the DRV-4.3 "PING
# test" row is cleared out).
# ---------------------------------------------------------------------------

# Row 7 - PS-1.6 now assigned, started (0.5 MD used) and In Progress
$ws.Range("A7").Value = 1.0
$ws.Range("E7").Value = "Anusree"
$ws.Range("G7").Value = 0.5
$ws.Range("I7").Value = "In Progress"

# Row 8 - PS-1.7, Sprint cell now populated (still blank value, just style)
$ws.Range("A8").Value = ""

# Row 9 - new User Story row: PS-1.8 / Integrate FreeRTOS into the project.
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = "PS-1.8"
$ws.Range("C9").Value = "PS-1"
$ws.Range("D9").Value = "Integrate FreeRTOS into the project."

# Row 10 - new User Story row: PS-1.9 / Create initialization task ...
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = "PS-1.9"
$ws.Range("C10").Value = "PS-1"
$ws.Range("D10").Value = "Create initialization task and start the scheduler."

# Row 11 - new User Story row: PS-1.10 / Debug print framework
$ws.Range("A11").Value = 1.0
$ws.Range("B11").Value = "PS-1.10"
$ws.Range("C11").Value = "PS-1"
$ws.Range("D11").Value = "Debug print framework"
$ws.Range("I11").Value = "To Do"

# Row 12 - DRV-2.1
$ws.Range("B12").Value = "DRV-2.1"
$ws.Range("C12").Value = "DRV-2"
$ws.Range("D12").Value = "Implement sensor power-up and basic I2C initialization sequence."

# Row 13 - DRV-2.2, now Done, Remaining MD 0.0
$ws.Range("B13").Value = "DRV-2.2"
$ws.Range("C13").Value = "DRV-2"
$ws.Range("D13").Value = "Configure sensor to operate in the required Gesture Mode."
$ws.Range("H13").Value = 0.0

# Row 14 - DRV-2.3, Sprint set
$ws.Range("A14").Value = 1.0
$ws.Range("B14").Value = "DRV-2.3"
$ws.Range("C14").Value = "DRV-2"
$ws.Range("D14").Value = "Implement read_gesture_id() function with stable state & Debug print."
$ws.Range("J14").Value = "Reads the classified integer."

# Row 15 - DRV-3.2, now Done
$ws.Range("B15").Value = "DRV-3.2"
$ws.Range("C15").Value = "DRV-3"
$ws.Range("D15").Value = "Implement gesture_feedback(state gesture) using the utility."
$ws.Range("I15").Value = "Done"
$ws.Range("J15").Value = "Controls the 5-LED output."

# Row 16 - DRV-4.1, now assigned to Jyothish, High priority, started
$ws.Range("A16").Value = 1.0
$ws.Range("B16").Value = "DRV-4.1"
$ws.Range("C16").Value = "DRV-4"
$ws.Range("D16").Value = "Integrate LWIP library with the Ethernet MAC HAL."
$ws.Range("E16").Value = "Jyothish"
$ws.Range("F16").Value = "High"
$ws.Range("G16").Value = 1.5

# Row 17 - DRV-4.2, now assigned to Jyothish, High priority, started
$ws.Range("A17").Value = 1.0
$ws.Range("B17").Value = "DRV-4.2"
$ws.Range("C17").Value = "DRV-4"
$ws.Range("D17").Value = "Configure LWIP for DHCP address acquisition."
$ws.Range("E17").Value = "Jyothish"
$ws.Range("F17").Value = "High"
$ws.Range("G17").Value = 1.0

# Row 18 - DRV-4.3 (PING test) task removed entirely
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""

# ---------------------------------------------------------------------------
# Formatting: reuse existing cell formats from equivalent cells via
# copy/paste-special (formats only) so the shared style entries line up with
# the ones already used elsewhere in the sheet.
# (NB: multi-area ranges such as "A7,A8,A9" only paste into the first area in
# this host, so each destination gets its own Copy/PasteSpecial pair.)
# ---------------------------------------------------------------------------

$formats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Style 8 (Sprint-style, font only) -> new Sprint cells + Assignee cells
foreach ($dest in @("A7","A8","A9","A10","A11","A16","A17","E7","E16","E17")) {
    $ws.Range("A6").Copy()
    $ws.Range($dest).PasteSpecial($formats)
}

# Style 3 (bordered input cell) -> new Task ID / User Story ID / Priority / Remaining MD / Estimated MD cells
foreach ($dest in @("B9","C9","B10","C10","B11","C11","F16","F17","G7","G16","G17","H13")) {
    $ws.Range("B2").Copy()
    $ws.Range($dest).PasteSpecial($formats)
}

# Style 9 (green fill, "In Progress") -> I7
$ws.Range("I14").Copy()
$ws.Range("I7").PasteSpecial($formats)

# Style 7 (red fill, "Done") -> I15
$ws.Range("I2").Copy()
$ws.Range("I15").PasteSpecial($formats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Data validation on the Status column now spans through row 17 (was row 15)
# ---------------------------------------------------------------------------
$v = $ws.Range("I2:I17").Validation
$v.Delete()
$v.Add(3, 1, 1, "=`$N`$3:`$N`$5")
$v.IgnoreBlank = $true
$v.ShowInput = $false
$v.ShowError = $true
